# rev 30/10 - ajustes graficos
# Applies the column-width / data-type tweaks made to the chart config
# workbook (g4_comparativo_faturamento_mes_combo.xlsx).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Configurações": narrow the "Título" column and switch the
# chart's value format from "currency" to "number".
# ---------------------------------------------------------------
$wsConfig = $wb.Worksheets.Item("Configurações")
$wsConfig.Columns.Item(2).ColumnWidth = 40
$wsConfig.Range("D2").Value = "number"

# ---------------------------------------------------------------
# Sheet "Dados": narrow the "Dataset" column and store every monthly
# figure as text instead of a number (format them as Text first so
# Excel doesn't silently convert the digits back into numbers).
# ---------------------------------------------------------------
$wsDados = $wb.Worksheets.Item("Dados")
$wsDados.Columns.Item(2).ColumnWidth = 20

$dadosValues = @{
    "C2" = "2129865.09";  "D2" = "9457954.52";  "E2" = "11211936.75"; "F2" = "8262338.15";
    "G2" = "11304130.34"; "H2" = "15653075.86"; "I2" = "19253061.11"; "J2" = "24540754.91";
    "K2" = "23931674.1";  "L2" = "31469128.71"; "M2" = "27907799.72"; "N2" = "33302632.57";

    "C3" = "3895831.63";  "D3" = "4918031.37";  "E3" = "7205905.36";  "F3" = "4301091.22";
    "G3" = "4433164.41";  "H3" = "5623704.78";  "I3" = "5529368.69";  "J3" = "12715598.51";
    "K3" = "13929117.72"; "L3" = "15887361.31"; "M3" = "17326931.71"; "N3" = "20974981.94";

    "C4" = "2129865.09";  "D4" = "11587819.61"; "E4" = "22799756.36"; "F4" = "31062094.51";
    "G4" = "42366224.85"; "H4" = "58019300.71"; "I4" = "77272361.82"; "J4" = "101813116.73";
    "K4" = "125744790.84";"L4" = "157213919.55";"M4" = "185121719.27";"N4" = "218424351.84";

    "C5" = "3895831.63";  "E5" = "16019768.36"; "F5" = "20320859.58";
    "G5" = "24754023.99"; "H5" = "30377728.77"; "I5" = "35907097.46"; "J5" = "48622695.97";
    "K5" = "62551813.69"; "M5" = "95766106.71"; "N5" = "116741088.65"
}

foreach ($addr in $dadosValues.Keys) {
    $cell = $wsDados.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $dadosValues[$addr]
}

# ---------------------------------------------------------------
# Sheet "Cores": rework the column widths (B/C/D/F resized, new G
# column added at 15.83).
# ---------------------------------------------------------------
$wsCores = $wb.Worksheets.Item("Cores")
$wsCores.Columns.Item(2).ColumnWidth = 20
$wsCores.Columns.Item(3).ColumnWidth = 20
$wsCores.Columns.Item(4).ColumnWidth = 15
$wsCores.Columns.Item(6).ColumnWidth = 12
$wsCores.Columns.Item(7).ColumnWidth = 15

# ---------------------------------------------------------------
# Sheet "Eixos": the manual min/max axis values are no longer used -
# clear them out.
# ---------------------------------------------------------------
$wsEixos = $wb.Worksheets.Item("Eixos")
$wsEixos.Range("C2:D3").ClearContents()
